$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value2 = "25.990.14"
$ws.Range("E2").Value2 = "  -0.35%  "
$ws.Range("D3").Value2 = "1.643.95"
$ws.Range("E3").Value2 = "  -0.29%  "
$ws.Range("D4").NumberFormat = "@"
$ws.Range("D4").Value2 = "1.005"
$ws.Range("D4").Style = $ws.Range("B4").Style
$ws.Range("E4").Value2 = "  -0.31%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value2 = "215.66"
$ws.Range("D5").Style = $ws.Range("B5").Style
$ws.Range("E5").Value2 = "  +0.14%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value2 = "0.5063"
$ws.Range("D6").Style = $ws.Range("B6").Style
$ws.Range("E6").Value2 = "  -0.68%  "
$ws.Range("E7").Value2 = "  -0.69%  "
$ws.Range("E8").Value2 = "  -0.50%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value2 = "0.06420"
$ws.Range("D9").Style = $ws.Range("B9").Style
$ws.Range("E9").Value2 = "  -0.08%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value2 = "19.68"
$ws.Range("D10").Style = $ws.Range("B10").Style
$ws.Range("E10").Value2 = "  +0.35%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value2 = "0.07755"
$ws.Range("D11").Style = $ws.Range("B11").Style
$ws.Range("E11").Value2 = "  +0.45%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value2 = "4.278"
$ws.Range("D12").Style = $ws.Range("B12").Style
$ws.Range("E12").Value2 = "  +0.09%  "
$ws.Range("D13").Value2 = "1.642.93"
$ws.Range("E13").Value2 = "  -0.52%  "
$ws.Range("D14").Value2 = "1.870.06"
$ws.Range("E14").Value2 = "  -0.40%  "
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value2 = "0.5467"
$ws.Range("D15").Style = $ws.Range("B15").Style
$ws.Range("E15").Value2 = "  +0.04%  "
$ws.Range("D16").Value2 = "0.0₅7951"
$ws.Range("E16").Value2 = "  -0.38%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value2 = "64.48"
$ws.Range("D17").Style = $ws.Range("B17").Style
$ws.Range("E17").Value2 = "  +1.04%  "
$ws.Range("D18").Value2 = "25.995.05"
$ws.Range("E18").Value2 = "  -0.41%  "
$ws.Range("E19").Value2 = "  -0.64%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value2 = "202.07"
$ws.Range("D20").Style = $ws.Range("B20").Style
$ws.Range("E20").Value2 = "  -2.46%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value2 = "4.393"
$ws.Range("D21").Style = $ws.Range("B21").Style
$ws.Range("E21").Value2 = "  +0.28%  "
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value2 = "9.906"
$ws.Range("D22").Style = $ws.Range("B22").Style
$ws.Range("E22").Value2 = "  -1.23%  "
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value2 = "5.995"
$ws.Range("D23").Style = $ws.Range("B23").Style
$ws.Range("E23").Value2 = "  -0.09%  "
$ws.Range("E24").Value2 = "  -0.66%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value2 = "1.872"
$ws.Range("D25").Style = $ws.Range("B25").Style
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value2 = "140.92"
$ws.Range("D26").Style = $ws.Range("B26").Style
$ws.Range("E26").Value2 = "  -1.54%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value2 = "0.1140"
$ws.Range("D27").Style = $ws.Range("B27").Style
$ws.Range("E27").Value2 = "  -2.27%  "
$ws.Range("B28").Value2 = "Cosmos"
$ws.Range("C28").Value2 = "https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom"
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value2 = "6.831"
$ws.Range("D28").Style = $ws.Range("B28").Style
$ws.Range("E28").Value2 = "  -1.07%  "
$ws.Range("B29").Value2 = "EthereumClassic"
$ws.Range("C29").Value2 = "https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc"
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value2 = "15.68"
$ws.Range("D29").Style = $ws.Range("B29").Style
$ws.Range("E29").Value2 = "  -0.74%  "
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value2 = "1.242"
$ws.Range("D30").Style = $ws.Range("B30").Style
$ws.Range("E30").Value2 = "  -0.08%  "
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value2 = "0.04934"
$ws.Range("D31").Style = $ws.Range("B31").Style
$ws.Range("E31").Value2 = "  -2.60%  "
$ws.Range("E32").Value2 = "  -1.69%  "
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value2 = "3.215"
$ws.Range("D33").Style = $ws.Range("B33").Style
$ws.Range("E33").Value2 = "  -0.34%  "
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value2 = "1.544"
$ws.Range("D34").Style = $ws.Range("B34").Style
$ws.Range("E34").Value2 = "  -0.22%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value2 = "2.365"
$ws.Range("D35").Style = $ws.Range("B35").Style
$ws.Range("E35").Value2 = "  +0.76%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value2 = "0.8946"
$ws.Range("D36").Style = $ws.Range("B36").Style
$ws.Range("E36").Value2 = "  -2.06%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value2 = "2.620"
$ws.Range("D37").Style = $ws.Range("B37").Style
$ws.Range("E37").Value2 = "  -1.03%  "
$ws.Range("D38").Value2 = "1.156.19"
$ws.Range("E38").Value2 = "  +1.12%  "
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value2 = "0.5593"
$ws.Range("D39").Style = $ws.Range("B39").Style
$ws.Range("E39").Value2 = "  -1.91%  "
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value2 = "0.01568"
$ws.Range("D40").Style = $ws.Range("B40").Style
$ws.Range("E40").Value2 = "  -0.37%  "
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value2 = "1.001"
$ws.Range("D41").Style = $ws.Range("B41").Style
$ws.Range("E41").Value2 = "  -0.65%  "
$ws.Range("E42").Value2 = "  +0.85%  "
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value2 = "0.8091"
$ws.Range("D43").Style = $ws.Range("B43").Style
$ws.Range("E43").Value2 = "  -1.68%  "
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value2 = "99.76"
$ws.Range("D44").Style = $ws.Range("B44").Style
$ws.Range("E44").Value2 = "  -0.12%  "
$ws.Range("D45").Value2 = "1.780.90"
$ws.Range("E45").Value2 = "  -0.45%  "
$ws.Range("E46").Value2 = "  +4.25%  "
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value2 = "0.4520"
$ws.Range("D47").Style = $ws.Range("B47").Style
$ws.Range("E47").Value2 = "  -0.36%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value2 = "1.002"
$ws.Range("D48").Style = $ws.Range("B48").Style
$ws.Range("E48").Value2 = "  -0.85%  "
$ws.Range("E49").Value2 = "  -0.69%  "
$ws.Range("D50").NumberFormat = "@"
$ws.Range("D50").Value2 = "0.05045"
$ws.Range("D50").Style = $ws.Range("B50").Style
$ws.Range("E50").Value2 = "  -0.55%  "
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value2 = "1.003"
$ws.Range("D51").Style = $ws.Range("B51").Style
$ws.Range("E51").Value2 = "  -0.55%  "
